# Revisi laporan BA SO Tahap 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SEMESTER : " -> "SEMESTER : [c.semester]"
$ws.Range("A7").Value = "SEMESTER : [c.semester]"

# " [a.subtotal_saldo_awal; block=row; when [a.cetak_subtotal]=1]" -> " [a.counter; block=row; when [a.cetak_subtotal]=1]"
$ws.Range("A14").Value = " [a.counter; block=row; when [a.cetak_subtotal]=1]"

# Split the combined "An Pengguna / Kuasa Pengguna Barang Pengurus Barang Pembantu Pengurus Barang"
# across A18 / F18 / L18 into three separate labels.
$ws.Range("A18").Value = "An Pengguna / Kuasa Pengguna Barang "
$ws.Range("F18").Value = "Pengurus Barang "
$ws.Range("L18").Value = "Pembantu Pengurus Barang"

# Column width tweaks (B, C, D:O all got nudged during this revision).
$ws.Columns.Item(2).ColumnWidth = 32.714285714285715
$ws.Columns.Item(3).ColumnWidth = 21.142857142857142
$ws.Range("D1:O1").EntireColumn.ColumnWidth = 17.714285714285715

# Sheet view: shift the top-left cell and the active selection.
$ws.Range("A7").Select()
